$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 13:52"

# Update updated case counts (Casos totales / Casos activos / Recuperados / Muertes)
$ws.Range("B6").Value = 6331
$ws.Range("C6").Value = 688
$ws.Range("D6").Value = 5385
$ws.Range("E6").Value = 258
$ws.Range("B20").Value = 1807
$ws.Range("C20").Value = 439
$ws.Range("D20").Value = 1151
$ws.Range("E20").Value = 217
$ws.Range("B23").Value = 1602
$ws.Range("C23").Value = 553
$ws.Range("D23").Value = 896
$ws.Range("E23").Value = 153
$ws.Range("B24").Value = 1560
$ws.Range("C24").Value = 3728
$ws.Range("D24").Value = 4707
$ws.Range("E24").Value = 82
$ws.Range("B25").Value = 1536
$ws.Range("C25").Value = 333
$ws.Range("D25").Value = 1411
$ws.Range("E25").Value = 30
$ws.Range("D26").Value = 1288
$ws.Range("E26").Value = 84
$ws.Range("B27").Value = 1501
$ws.Range("C27").Value = 129
$ws.Range("D27").Value = 1258
$ws.Range("E27").Value = 114
$ws.Range("B29").Value = 1349
$ws.Range("C29").Value = 387
$ws.Range("D29").Value = 847
$ws.Range("E29").Value = 115
$ws.Range("B30").Value = 1344
$ws.Range("C30").Value = 564
$ws.Range("D30").Value = 584
$ws.Range("E30").Value = 196
$ws.Range("B31").Value = 1283
$ws.Range("C31").Value = 193
$ws.Range("D31").Value = 1012
$ws.Range("E31").Value = 78
$ws.Range("B34").Value = 1024
$ws.Range("C34").Value = 387
$ws.Range("D34").Value = 517
$ws.Range("E34").Value = 120
$ws.Range("B40").Value = 837
$ws.Range("C40").Value = 172
$ws.Range("D40").Value = 601
$ws.Range("B41").Value = 766
$ws.Range("C41").Value = 253
$ws.Range("D41").Value = 430
$ws.Range("E41").Value = 83
$ws.Range("B42").Value = 751
$ws.Range("C42").Value = 333
$ws.Range("D42").Value = 660
$ws.Range("E42").Value = 22
$ws.Range("B46").Value = 510
$ws.Range("C46").Value = 126
$ws.Range("D46").Value = 343
$ws.Range("E46").Value = 41
$ws.Range("B51").Value = 342
$ws.Range("C51").Value = 107
$ws.Range("D51").Value = 196
$ws.Range("E51").Value = 39
